# Natmi following Dr Hou advice
# Update LR-pair stats (Psen1-Notch4) for rows 2-17: expressing-cell counts
# (columns E and K) move from 1 to 3 cells, and the dependent expression /
# specificity statistics are recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.506869
$ws.Range("H2").Value = 49.520607
$ws.Range("I2").Value = 0.2165594803671733
$ws.Range("J2").Value = 0.2165594803671733
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 46.05975733333333
$ws.Range("N2").Value = 138.179272
$ws.Range("O2").Value = 0.8970651351272991
$ws.Range("P2").Value = 0.897065135127299
$ws.Range("Q2").Value = 760.3023804731225
$ws.Range("R2").Value = 6842.721424258104
$ws.Range("S2").Value = 0.194267959518676
$ws.Range("T2").Value = 0.194267959518676

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.506869
$ws.Range("H3").Value = 49.520607
$ws.Range("I3").Value = 0.2165594803671733
$ws.Range("J3").Value = 0.2165594803671733
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.770761666666667
$ws.Range("N3").Value = 8.312284999999999
$ws.Range("O3").Value = 0.05396367312415441
$ws.Range("P3").Value = 0.0539636731241544
$ws.Range("Q3").Value = 45.73659986188833
$ws.Range("R3").Value = 411.6293987569949
$ws.Range("S3").Value = 0.01168634501047087
$ws.Range("T3").Value = 0.01168634501047087

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.506869
$ws.Range("H4").Value = 49.520607
$ws.Range("I4").Value = 0.2165594803671733
$ws.Range("J4").Value = 0.2165594803671733
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.802173666666667
$ws.Range("N4").Value = 5.406521
$ws.Range("O4").Value = 0.03509934175535083
$ws.Range("P4").Value = 0.03509934175535083
$ws.Range("Q4").Value = 29.74824463091633
$ws.Range("R4").Value = 267.734201678247
$ws.Range("S4").Value = 0.007601095211768604
$ws.Range("T4").Value = 0.007601095211768605

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.506869
$ws.Range("H5").Value = 49.520607
$ws.Range("I5").Value = 0.2165594803671733
$ws.Range("J5").Value = 0.2165594803671733
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7122493333333333
$ws.Range("N5").Value = 2.136748
$ws.Range("O5").Value = 0.01387184999319569
$ws.Range("P5").Value = 0.01387184999319569
$ws.Range("Q5").Value = 11.75700644067066
$ws.Range("R5").Value = 105.813057966036
$ws.Range("S5").Value = 0.003004080626257836
$ws.Range("T5").Value = 0.003004080626257836

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 24.781512
$ws.Range("H6").Value = 74.34453600000001
$ws.Range("I6").Value = 0.3251174623990092
$ws.Range("J6").Value = 0.3251174623990092
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 46.05975733333333
$ws.Range("N6").Value = 138.179272
$ws.Range("O6").Value = 0.8970651351272991
$ws.Range("P6").Value = 0.897065135127299
$ws.Range("Q6").Value = 1141.430429073088
$ws.Range("R6").Value = 10272.87386165779
$ws.Range("S6").Value = 0.2916515403392118
$ws.Range("T6").Value = 0.2916515403392118

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 24.781512
$ws.Range("H7").Value = 74.34453600000001
$ws.Range("I7").Value = 0.3251174623990092
$ws.Range("J7").Value = 0.3251174623990092
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.770761666666667
$ws.Range("N7").Value = 8.312284999999999
$ws.Range("O7").Value = 0.05396367312415441
$ws.Range("P7").Value = 0.0539636731241544
$ws.Range("Q7").Value = 68.66366349164001
$ws.Range("R7").Value = 617.97297142476
$ws.Range("S7").Value = 0.01754453246785469
$ws.Range("T7").Value = 0.01754453246785469

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 24.781512
$ws.Range("H8").Value = 74.34453600000001
$ws.Range("I8").Value = 0.3251174623990092
$ws.Range("J8").Value = 0.3251174623990092
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.802173666666667
$ws.Range("N8").Value = 5.406521
$ws.Range("O8").Value = 0.03509934175535083
$ws.Range("P8").Value = 0.03509934175535083
$ws.Range("Q8").Value = 44.660588346584
$ws.Range("R8").Value = 401.945295119256
$ws.Range("S8").Value = 0.01141140892337525
$ws.Range("T8").Value = 0.01141140892337525

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 24.781512
$ws.Range("H9").Value = 74.34453600000001
$ws.Range("I9").Value = 0.3251174623990092
$ws.Range("J9").Value = 0.3251174623990092
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.7122493333333333
$ws.Range("N9").Value = 2.136748
$ws.Range("O9").Value = 0.01387184999319569
$ws.Range("P9").Value = 0.01387184999319569
$ws.Range("Q9").Value = 17.650615400992
$ws.Range("R9").Value = 158.855538608928
$ws.Range("S9").Value = 0.004509980668567497
$ws.Range("T9").Value = 0.004509980668567497

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.67943933333333
$ws.Range("H10").Value = 74.038318
$ws.Range("I10").Value = 0.3237783348120013
$ws.Range("J10").Value = 0.3237783348120013
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.05975733333333
$ws.Range("N10").Value = 138.179272
$ws.Range("O10").Value = 0.8970651351272991
$ws.Range("P10").Value = 0.897065135127299
$ws.Range("Q10").Value = 1136.728986816055
$ws.Range("R10").Value = 10230.5608813445
$ws.Range("S10").Value = 0.2904502556694198
$ws.Range("T10").Value = 0.2904502556694198

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 24.67943933333333
$ws.Range("H11").Value = 74.038318
$ws.Range("I11").Value = 0.3237783348120013
$ws.Range("J11").Value = 0.3237783348120013
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.770761666666667
$ws.Range("N11").Value = 8.312284999999999
$ws.Range("O11").Value = 0.05396367312415441
$ws.Range("P11").Value = 0.0539636731241544
$ws.Range("Q11").Value = 68.38084445962555
$ws.Range("R11").Value = 615.42760013663
$ws.Range("S11").Value = 0.01747226822447786
$ws.Range("T11").Value = 0.01747226822447786

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 24.67943933333333
$ws.Range("H12").Value = 74.038318
$ws.Range("I12").Value = 0.3237783348120013
$ws.Range("J12").Value = 0.3237783348120013
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.802173666666667
$ws.Range("N12").Value = 5.406521
$ws.Range("O12").Value = 0.03509934175535083
$ws.Range("P12").Value = 0.03509934175535083
$ws.Range("Q12").Value = 44.47663567463089
$ws.Range("R12").Value = 400.289721071678
$ws.Range("S12").Value = 0.01136440642654484
$ws.Range("T12").Value = 0.01136440642654484

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 24.67943933333333
$ws.Range("H13").Value = 74.038318
$ws.Range("I13").Value = 0.3237783348120013
$ws.Range("J13").Value = 0.3237783348120013
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7122493333333333
$ws.Range("N13").Value = 2.136748
$ws.Range("O13").Value = 0.01387184999319569
$ws.Range("P13").Value = 0.01387184999319569
$ws.Range("Q13").Value = 17.57791421220711
$ws.Range("R13").Value = 158.201227909864
$ws.Range("S13").Value = 0.004491404491558774
$ws.Range("T13").Value = 0.004491404491558773

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.255437
$ws.Range("H14").Value = 30.766311
$ws.Range("I14").Value = 0.1345447224218162
$ws.Range("J14").Value = 0.1345447224218162
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 46.05975733333333
$ws.Range("N14").Value = 138.179272
$ws.Range("O14").Value = 0.8970651351272991
$ws.Range("P14").Value = 0.897065135127299
$ws.Range("Q14").Value = 472.362939567288
$ws.Range("R14").Value = 4251.266456105592
$ws.Range("S14").Value = 0.1206953795999915
$ws.Range("T14").Value = 0.1206953795999915

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.255437
$ws.Range("H15").Value = 30.766311
$ws.Range("I15").Value = 0.1345447224218162
$ws.Range("J15").Value = 0.1345447224218162
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.770761666666667
$ws.Range("N15").Value = 8.312284999999999
$ws.Range("O15").Value = 0.05396367312415441
$ws.Range("P15").Value = 0.0539636731241544
$ws.Range("Q15").Value = 28.415371714515
$ws.Range("R15").Value = 255.738345430635
$ws.Range("S15").Value = 0.007260527421350978
$ws.Range("T15").Value = 0.007260527421350978

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.255437
$ws.Range("H16").Value = 30.766311
$ws.Range("I16").Value = 0.1345447224218162
$ws.Range("J16").Value = 0.1345447224218162
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.802173666666667
$ws.Range("N16").Value = 5.406521
$ws.Range("O16").Value = 0.03509934175535083
$ws.Range("P16").Value = 0.03509934175535083
$ws.Range("Q16").Value = 18.482078501559
$ws.Range("R16").Value = 166.338706514031
$ws.Range("S16").Value = 0.00472243119366214
$ws.Range("T16").Value = 0.004722431193662141

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.255437
$ws.Range("H17").Value = 30.766311
$ws.Range("I17").Value = 0.1345447224218162
$ws.Range("J17").Value = 0.1345447224218162
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.7122493333333333
$ws.Range("N17").Value = 2.136748
$ws.Range("O17").Value = 0.01387184999319569
$ws.Range("P17").Value = 0.01387184999319569
$ws.Range("Q17").Value = 7.304428166292
$ws.Range("R17").Value = 65.739853496628
$ws.Range("S17").Value = 0.001866384206811588
$ws.Range("T17").Value = 0.001866384206811588
